$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "floodmedia" column (H) was left as numeric 0 placeholders for every
# data row; replace them with the text "None" to reflect that no flood
# media was actually used.
$ws.Range("H2:H37").Value = "None"

# Row 2 had a stale explicit row height (16) left over from before; once its
# content is touched it settles back down to the same height used by every
# other data row (15).
$ws.Rows.Item(2).RowHeight = 15

# Move the viewport/selection so the newly edited floodmedia column is in
# view and selected, matching where the author was working.
$excel.Goto($ws.Range("A8"), $false) | Out-Null
$ws.Range("H2:H37").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
